# Scheduled-runner style update of market-price derived columns
# (currentAveragePrice*, LevePrice*, LeveProfit*) across several
# worksheets, refreshing them with newly-fetched values.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H6").Value = 98.59999999999999
$ws.Range("I6").Value = 98.59999999999999
$ws.Range("K6").Value = 295.8
$ws.Range("M6").Value = -183.8

$ws.Range("H43").Value = 3754.2727
$ws.Range("I43").Value = 2750
$ws.Range("J43").Value = 3977.4443
$ws.Range("K43").Value = 2750
$ws.Range("L43").Value = 3977.4443
$ws.Range("M43").Value = -2681
$ws.Range("N43").Value = -4115.4443

$ws.Range("H64").Value = 4711.6
$ws.Range("J64").Value = 5000
$ws.Range("L64").Value = 5000
$ws.Range("N64").Value = -5496

$ws.Range("H67").Value = 4711.6
$ws.Range("J67").Value = 5000
$ws.Range("L67").Value = 5000
$ws.Range("N67").Value = -6716

$ws.Range("H98").Value = 2128
$ws.Range("I98").Value = 2128
$ws.Range("J98").Value = 0
$ws.Range("K98").Value = 2128
$ws.Range("L98").Value = 0
$ws.Range("M98").Value = -630
$ws.Range("N98").ClearContents()

$ws.Range("H113").Value = 4601.125
$ws.Range("J113").Value = 4860.8
$ws.Range("L113").Value = 4860.8
$ws.Range("N113").Value = -11368.8

$ws.Range("H122").Value = 2128
$ws.Range("I122").Value = 2128
$ws.Range("J122").Value = 0
$ws.Range("K122").Value = 6384
$ws.Range("L122").Value = 0
$ws.Range("M122").Value = -3934
$ws.Range("N122").ClearContents()

$ws.Range("H138").Value = 4870.9697
$ws.Range("I138").Value = 3665.077
$ws.Range("J138").Value = 5654.8
$ws.Range("K138").Value = 10995.231
$ws.Range("L138").Value = 16964.4
$ws.Range("M138").Value = -5855.231
$ws.Range("N138").Value = -27244.4

$ws.Range("H141").Value = 3248.9583
$ws.Range("I141").Value = 1474.2222
$ws.Range("J141").Value = 8573.166999999999
$ws.Range("K141").Value = 4422.6666
$ws.Range("L141").Value = 25719.501
$ws.Range("M141").Value = 757.3334000000004
$ws.Range("N141").Value = -36079.501

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H7").Value = 72500
$ws.Range("J7").Value = 90000
$ws.Range("L7").Value = 90000
$ws.Range("N7").Value = -90228

$ws.Range("H45").Value = 6079.231
$ws.Range("I45").Value = 26953
$ws.Range("J45").Value = 2284
$ws.Range("K45").Value = 26953
$ws.Range("L45").Value = 2284
$ws.Range("M45").Value = -26576
$ws.Range("N45").Value = -3038

$ws.Range("H61").Value = 54692.95
$ws.Range("I61").Value = 1827.8
$ws.Range("J61").Value = 252937.25
$ws.Range("K61").Value = 1827.8
$ws.Range("L61").Value = 252937.25
$ws.Range("M61").Value = -1615.8
$ws.Range("N61").Value = -253361.25

$ws.Range("H122").Value = 3129.6296
$ws.Range("I122").Value = 2615.875
$ws.Range("J122").Value = 3876.9092
$ws.Range("K122").Value = 7847.625
$ws.Range("L122").Value = 11630.7276
$ws.Range("M122").Value = -5397.625
$ws.Range("N122").Value = -16530.7276

$ws.Range("H136").Value = 54692.95
$ws.Range("I136").Value = 1827.8
$ws.Range("J136").Value = 252937.25
$ws.Range("K136").Value = 5483.4
$ws.Range("L136").Value = 758811.75
$ws.Range("M136").Value = -2933.4
$ws.Range("N136").Value = -763911.75

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H64").Value = 112364
$ws.Range("J64").Value = 1817.2
$ws.Range("L64").Value = 1817.2
$ws.Range("N64").Value = -2267.2

$ws.Range("H67").Value = 112364
$ws.Range("J67").Value = 1817.2
$ws.Range("L67").Value = 1817.2
$ws.Range("N67").Value = -3377.2

$ws.Range("H99").Value = 1252811.2
$ws.Range("I99").Value = 2358.3333
$ws.Range("J99").Value = 2407075.5
$ws.Range("K99").Value = 2358.3333
$ws.Range("L99").Value = 2407075.5
$ws.Range("M99").Value = -860.3332999999998
$ws.Range("N99").Value = -2410071.5

$ws.Range("H105").Value = 2548.4285
$ws.Range("I105").Value = 2147.5
$ws.Range("K105").Value = 2147.5
$ws.Range("M105").Value = -400.5

$ws.Range("H134").Value = 2407.6667
$ws.Range("I134").Value = 2213.7368
$ws.Range("K134").Value = 6641.2104
$ws.Range("M134").Value = -4106.2104

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 1545.4286
$ws.Range("I16").Value = 776.25
$ws.Range("J16").Value = 2571
$ws.Range("K16").Value = 776.25
$ws.Range("L16").Value = 2571
$ws.Range("M16").Value = -489.25
$ws.Range("N16").Value = -3145

$ws.Range("H31").Value = 4690.811
$ws.Range("I31").Value = 2658.6978
$ws.Range("J31").Value = 6549.9785
$ws.Range("K31").Value = 2658.6978
$ws.Range("L31").Value = 6549.9785
$ws.Range("M31").Value = -2363.6978
$ws.Range("N31").Value = -7139.9785

$ws.Range("H34").Value = 4690.811
$ws.Range("I34").Value = 2658.6978
$ws.Range("J34").Value = 6549.9785
$ws.Range("K34").Value = 2658.6978
$ws.Range("L34").Value = 6549.9785
$ws.Range("M34").Value = -2456.6978
$ws.Range("N34").Value = -6953.9785

$ws.Range("H58").Value = 3249.5
$ws.Range("I58").Value = 3249.5
$ws.Range("J58").Value = 0
$ws.Range("K58").Value = 3249.5
$ws.Range("L58").Value = 0
$ws.Range("M58").Value = -3046.5
$ws.Range("N58").ClearContents()

$ws.Range("H99").Value = 4467493
$ws.Range("J99").Value = 4467493
$ws.Range("L99").Value = 4467493
$ws.Range("N99").Value = -4470489

$ws.Range("H107").Value = 945.2
$ws.Range("I107").Value = 826.931
$ws.Range("J107").Value = 1516.8334
$ws.Range("K107").Value = 826.931
$ws.Range("L107").Value = 1516.8334
$ws.Range("M107").Value = 1093.069
$ws.Range("N107").Value = -5356.8334

$ws.Range("H113").Value = 1545.4286
$ws.Range("I113").Value = 776.25
$ws.Range("J113").Value = 2571
$ws.Range("K113").Value = 776.25
$ws.Range("L113").Value = 2571
$ws.Range("M113").Value = 1393.75
$ws.Range("N113").Value = -6911

$ws.Range("H126").Value = 4467493
$ws.Range("J126").Value = 4467493
$ws.Range("L126").Value = 13402479
$ws.Range("N126").Value = -13407419

$ws.Range("H132").Value = 2456
$ws.Range("I132").Value = 2129.6155
$ws.Range("K132").Value = 6388.8465
$ws.Range("M132").Value = -3858.8465

$ws.Range("H136").Value = 3249.5
$ws.Range("I136").Value = 3249.5
$ws.Range("J136").Value = 0
$ws.Range("K136").Value = 9748.5
$ws.Range("L136").Value = 0
$ws.Range("M136").Value = -7198.5
$ws.Range("N136").ClearContents()

$ws.Range("H141").Value = 323288.9
$ws.Range("J141").Value = 323288.9
$ws.Range("L141").Value = 323288.9
$ws.Range("N141").Value = -333648.9

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 2099.3
$ws.Range("I5").Value = 0
$ws.Range("J5").Value = 2099.3
$ws.Range("K5").Value = 0
$ws.Range("L5").Value = 6297.900000000001
$ws.Range("M5").ClearContents()
$ws.Range("N5").Value = -6521.900000000001

$ws.Range("H55").Value = 142867280
$ws.Range("J55").Value = 142867280
$ws.Range("L55").Value = 428601840
$ws.Range("N55").Value = -428602194

$ws.Range("H131").Value = 113176.78
$ws.Range("J131").Value = 3168
$ws.Range("L131").Value = 9504
$ws.Range("N131").Value = -19584

$ws.Range("H135").Value = 2099.3
$ws.Range("I135").Value = 0
$ws.Range("J135").Value = 2099.3
$ws.Range("K135").Value = 0
$ws.Range("L135").Value = 18893.7
$ws.Range("M135").ClearContents()
$ws.Range("N135").Value = -23963.7

$ws.Range("H139").Value = 5435.0347
$ws.Range("I139").Value = 2267.3125
$ws.Range("K139").Value = 6801.9375
$ws.Range("M139").Value = -1661.9375

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 2624.8948
$ws.Range("I102").Value = 2158.1428
$ws.Range("J102").Value = 3931.8
$ws.Range("K102").Value = 2158.1428
$ws.Range("L102").Value = 3931.8
$ws.Range("M102").Value = -536.1428000000001
$ws.Range("N102").Value = -7175.8

$ws.Range("H122").Value = 7038.077
$ws.Range("I122").Value = 7560.5454
$ws.Range("K122").Value = 22681.6362
$ws.Range("M122").Value = -20231.6362

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H55").Value = 1221.8462
$ws.Range("J55").Value = 1922.3636
$ws.Range("L55").Value = 1922.3636
$ws.Range("N55").Value = -2268.3636

$ws.Range("H61").Value = 3387.5557
$ws.Range("I61").Value = 3198.2
$ws.Range("J61").Value = 3624.25
$ws.Range("K61").Value = 3198.2
$ws.Range("L61").Value = 3624.25
$ws.Range("M61").Value = -2996.2
$ws.Range("N61").Value = -4028.25

$ws.Range("H82").Value = 1298.5714
$ws.Range("I82").Value = 1201
$ws.Range("J82").Value = 1337.6
$ws.Range("K82").Value = 1201
$ws.Range("L82").Value = 1337.6
$ws.Range("M82").Value = -840
$ws.Range("N82").Value = -2059.6

$ws.Range("H85").Value = 1298.5714
$ws.Range("I85").Value = 1201
$ws.Range("J85").Value = 1337.6
$ws.Range("K85").Value = 1201
$ws.Range("L85").Value = 1337.6
$ws.Range("M85").Value = 47
$ws.Range("N85").Value = -3833.6

$ws.Range("H100").Value = 51400
$ws.Range("I100").Value = 51400
$ws.Range("J100").Value = 0
$ws.Range("K100").Value = 51400
$ws.Range("L100").Value = 0
$ws.Range("M100").Value = -50859
$ws.Range("N100").ClearContents()

$ws.Range("H113").Value = 3387.5557
$ws.Range("I113").Value = 3198.2
$ws.Range("J113").Value = 3624.25
$ws.Range("K113").Value = 3198.2
$ws.Range("L113").Value = 3624.25
$ws.Range("M113").Value = -1028.2
$ws.Range("N113").Value = -7964.25

$ws.Range("H122").Value = 16684.438
$ws.Range("I122").Value = 18094.727
$ws.Range("J122").Value = 13581.8
$ws.Range("K122").Value = 54284.181
$ws.Range("L122").Value = 40745.39999999999
$ws.Range("M122").Value = -51834.181
$ws.Range("N122").Value = -45645.39999999999

$ws.Range("H136").Value = 5696.3413
$ws.Range("I136").Value = 6627.048
$ws.Range("K136").Value = 19881.144
$ws.Range("M136").Value = -17331.144

$ws.Range("H139").Value = 66114.5
$ws.Range("I139").Value = 51285.715
$ws.Range("J139").Value = 100715
$ws.Range("K139").Value = 51285.715
$ws.Range("L139").Value = 100715
$ws.Range("M139").Value = -46145.715
$ws.Range("N139").Value = -110995

$ws.Range("H141").Value = 107290.6
$ws.Range("J141").Value = 107290.6
$ws.Range("L141").Value = 107290.6
$ws.Range("N141").Value = -117650.6
